$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.653.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.329.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.10%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'581.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.96%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'176.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.46%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.325.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.18%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.06%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.577"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.60%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'45.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.10%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -2.13%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'674.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +5.37%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.875.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.26%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.64%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.652.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.41%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.74%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.333.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.19%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.28%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.18%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +9.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'17.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'99.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.71%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.68%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'33.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.25%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.70%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +9.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'573.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.02%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.694.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'56.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.78%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'InjectiveProtocol"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'34.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Kaspa"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.52%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.15%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.65%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'PEPE"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.0₃0667"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.37%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'ApeXProtocol"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.34%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.38%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.65%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.85%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.41%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.35%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.51%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'129.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.18%  "
$ws.Range("E51").Style = "Normal"
